# 3.5.1.1.xlsx - add a new "2021" data column (M) mirroring the existing
# 2020 column (L), including number formats/styles, then restore the
# originally-saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clone the formatting of the whole 2020 column (L2:L33) onto the new
#    2021 column (M2:M33). This brings along number formats, fonts,
#    alignment and borders in one shot, matching the style index that
#    each row already uses in column L.
$ws.Range("L2:L33").Copy()
$ws.Range("M2:M33").PasteSpecial(-4122)

# 2) Header row: year label for the new column.
$ws.Range("M3").Value = 2021

# 3) Data values for the new column, row by row (values taken from the
#    2021 dataset; "-" denotes the existing "no data" marker already used
#    throughout the sheet).
$ws.Range("M4").Value = 2.0173148373954581
$ws.Range("M5").Value = 0.11867182493532386
$ws.Range("M6").Value = 3.9440914499323179
$ws.Range("M7").Value = 0
$ws.Range("M8").Value = "-"
$ws.Range("M9").Value = 0
$ws.Range("M10").Value = 0.62921030174566528
$ws.Range("M11").Value = "-"
$ws.Range("M12").Value = 1.2497227177719943
$ws.Range("M13").Value = 0.19844537890168421
$ws.Range("M14").Value = "-"
$ws.Range("M15").Value = 0.39861918314956984
$ws.Range("M16").Value = 0
$ws.Range("M17").Value = "-"
$ws.Range("M18").Value = 0
$ws.Range("M19").Value = 0.85521252031129735
$ws.Range("M20").Value = "-"
$ws.Range("M21").Value = 1.6913581464969858
$ws.Range("M22").Value = 1.8347815875998121
$ws.Range("M23").Value = "-"
$ws.Range("M24").Value = 3.6321107648498847
$ws.Range("M25").Value = 6.1211560415300026
$ws.Range("M26").Value = "-"
$ws.Range("M27").Value = 12.437939862560766
$ws.Range("M28").Value = 3.6823562661275693
$ws.Range("M29").Value = 0.69433233870225819
$ws.Range("M30").Value = 7.0564990356117976
$ws.Range("M31").Value = 2.7447727328177227
$ws.Range("M32").Value = "-"
$ws.Range("M33").Value = 5.6418550419377889

# 4) Two cells (M26 and M32) end up with a style different from their L
#    counterparts. Re-apply formatting from cells that already carry the
#    correct target style so the matching cellXfs entry is (re)used:
#      - M26 needs the numeric (0.0) style right-aligned -> derive it from
#        L26 (same numeric style) by nudging its alignment.
#      - M32 needs the plain "-" style already used by column L's other
#        "-" cells (e.g. L8), so copy that format directly.
$ws.Range("M26").HorizontalAlignment = -4152

$ws.Range("L8").Copy()
$ws.Range("M32").PasteSpecial(-4122)

# 5) Restore the selection that was active when the workbook was saved.
$ws.Range("P6").Select()
